# Applies the "adding new version of project report" edit to the ETL
# Project Report document:
#   - Paragraph "For the ETL Project...": splits/reworks the dataset
#     introduction sentences and the cleanup-reason sentence.
#   - Paragraph "Jupyter notebooks...": flattens the Jupyter sentence,
#     adds detail about dropping duplicates / the _GoBack bookmark moves
#     here.
#   - Paragraph "To load the data...": expands the troubleshooting
#     sentence and removes the (now relocated) _GoBack bookmark.

$d = $word.ActiveDocument

function Replace-ParagraphRuns($paragraphIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $r = $p.Range
    # Exclude the trailing paragraph mark so only the run content is
    # replaced; pPr (spacing/indent/etc.) on the paragraph stays intact.
    $target = $d.Range($r.Start, $r.End - 1)
    $pkg = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
           "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
           "<pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:body><w:p>" + $innerXml + "</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $target.InsertXML($pkg)
}

# --- Paragraph: "For the ETL Project, two CSV files were chosen..." ---
$para1 = '<w:r><w:t>For the ETL Project, two CSV files were chosen</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
         '<w:r><w:t>The first dataset</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> reflected data for crimes in </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">Los Angeles and was downloaded from Kaggle. The </w:t></w:r>' +
         '<w:r><w:t>second dataset</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> reflected data for crimes in </w:t></w:r>' +
         '<w:r><w:t>Kansas City and was downloaded from Open Data KC</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">. The data was not collected in the exact same manner in each city, so cleanup was needed </w:t></w:r>' +
         '<w:r><w:t>to keep the tables consistent with each other</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">. The idea was to collect this data </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">as </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">if a group were wanting to be able to compare </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">crime between </w:t></w:r>' +
         '<w:r><w:t>the two cities.</w:t></w:r>'

# --- Paragraph: "Jupyter notebooks was used..." ---
$para2 = '<w:r><w:t xml:space="preserve">Jupyter notebooks was used to read in the CSV files as pandas data frames. The </w:t></w:r>' +
         '<w:r><w:t>cleanup</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> involved changing column names since they had similar data but different labeling for each source. </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">We dropped duplicate entries </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">based on the incident ID column </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">and created a column </w:t></w:r>' +
         '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
         '<w:r><w:t xml:space="preserve">that reflected which crime data belonged to each individual city. MySQL was used to create </w:t></w:r>' +
         '<w:bookmarkEnd w:id="0"/>' +
         '<w:r><w:t xml:space="preserve">queries for the two different data files. </w:t></w:r>'

# --- Paragraph: "To load the data, we created a connection..." ---
$para3 = '<w:r><w:t xml:space="preserve">To load the data, we created a connection with the database by using the local host and then the create engine function. After confirming the two tables in MySQL, we loaded the </w:t></w:r>' +
         '<w:r><w:t>d</w:t></w:r>' +
         '<w:r><w:t>ata</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> fr</w:t></w:r>' +
         '<w:r><w:t>ames by using the &#8220;.</w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>to_sql</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t>&#8221; function.</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> We ran into a problem here. The connection would cause an error and drop our primary key. </w:t></w:r>' +
         '<w:r><w:t>We did some troubleshooting and found t</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">he reason </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">for the error </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">was because </w:t></w:r>' +
         '<w:r><w:t xml:space="preserve">we had the </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>incident_id</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> set as our index. Once we changed this, the connection filled through properly.</w:t></w:r>'

# Locate the three paragraphs by their distinctive leading text so the
# script is resilient to the exact paragraph index.
$targetIdx1 = 0
$targetIdx2 = 0
$targetIdx3 = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("For the ETL Project")) { $targetIdx1 = $i }
    if ($t.StartsWith("Jupyter notebooks")) { $targetIdx2 = $i }
    if ($t.StartsWith("To load the data")) { $targetIdx3 = $i }
}

Replace-ParagraphRuns $targetIdx1 $para1
Replace-ParagraphRuns $targetIdx2 $para2
Replace-ParagraphRuns $targetIdx3 $para3

Write-Output "Edit applied."
